$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Simple in-place text replacements (unique anchors, safe in any order)
# ------------------------------------------------------------------

# "--> [测试全部完成] "显示个人ID 和 匹配对象ID输入空格"" -> "--> [测试全部完成] 询问付费"
$d.Content.Find.Execute(
    '-->[测试全部完成] "显示个人ID 和 匹配对象ID输入空格"', $false, $false, $false, $false, $false,
    $true, 1, $false, '-->[测试全部完成] 询问付费', 2) | Out-Null

# "-->输入匹配对象ID并开始匹配" -> "-->付费后输入匹配对象ID并开始匹配"
$d.Content.Find.Execute(
    '-->输入匹配对象ID并开始匹配', $false, $false, $false, $false, $false,
    $true, 1, $false, '-->付费后输入匹配对象ID并开始匹配', 2) | Out-Null

# "-->检查用户是否已付费" -> "--> (*)"
$d.Content.Find.Execute(
    '-->检查用户是否已付费', $false, $false, $false, $false, $false,
    $true, 1, $false, '--> (*)', 2) | Out-Null

# "-->"显示个人ID 和 匹配对象ID输入空格"" (near the end) -> "-->"询问付费""
# (by now the only remaining occurrence of this exact fragment is the one near the end)
$d.Content.Find.Execute(
    '-->"显示个人ID 和 匹配对象ID输入空格"', $false, $false, $false, $false, $false,
    $true, 1, $false, '-->"询问付费"', 2) | Out-Null

# ------------------------------------------------------------------
# 2) Remove the old nested "双方用户均已付费 / 付费" branch block
#    (everything from "-->[双方均已付费]..." through the "endif" that
#    closes "if 双方用户均已付费 then", i.e. 20 whole paragraphs) and
#    collapse the "if 双方用户均已付费 then" paragraph itself down to
#    blank trailing whitespace.
# ------------------------------------------------------------------

$ifRange = $d.Content.Duplicate
$ifRange.Find.Execute('if 双方用户均已付费 then') | Out-Null
$ifPara = $ifRange.Paragraphs.Item(1)

$blockStartRange = $d.Content.Duplicate
$blockStartRange.Find.Execute('-->[双方均已付费]"显示配对结果并发送结果给导师"') | Out-Null
$blockStartIdx = $blockStartRange.Paragraphs.Item(1).Index
$blockEndIdx = $blockStartIdx + 19

$blockStartPara = $d.Paragraphs.Item($blockStartIdx)
$blockEndPara = $d.Paragraphs.Item($blockEndIdx)
$delRange = $d.Range($blockStartPara.Range.Start, $blockEndPara.Range.End)
$delRange.Delete() | Out-Null

# Now turn the former "if 双方用户均已付费 then" paragraph into a blank line
$ifPara.Range.Text = "       "

# ------------------------------------------------------------------
# 3) Insert the new "完成所有测试" line right after the
#    "-->[选择从历史进度开始] ..." paragraph.
# ------------------------------------------------------------------

$insRange = $d.Content.Duplicate
$insRange.Find.Execute('-->[选择从历史进度开始] “从历史进度开始”') | Out-Null
$insRange.Expand(4) | Out-Null
$insRange.Collapse(0) | Out-Null
$insRange.InsertBefore("          -->[完成所有测试] ""测试完成,显示测试结果"" `r")
